# Chef item edit: update the "Coffee" row (row 3) in the Item Data sheet.
# Price stays numerically 10, but is now written as an explicit numeric value,
# and the "Special Item" flag is switched on (False -> True).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 10.0
$ws.Range("D3").Value = $true
